$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 14 (shifts the old "Biodiesel", "Electricity", "Heat"
# rows down from 14-16 to 16-18), to add two new subfuels under "16_others":
# Municipal Solid Waste (Renewable) and Municipal Solid Waste (Non-renewable).
$ws.Rows("14:15").Insert()

# Populate column B ("fuels" code) first for both new rows - matches existing
# "16_others" value already present elsewhere in the sheet.
$ws.Cells.Item(14, 2).Value2 = "16_others"
$ws.Cells.Item(15, 2).Value2 = "16_others"

# Populate column C ("subfuels" code) next for both new rows.
$ws.Cells.Item(14, 3).Value2 = "16_03_municipal_solid_waste_renewable"
$ws.Cells.Item(15, 3).Value2 = "16_04_municipal_solid_waste_nonrenewable"

# Populate column A (display name) last for both new rows.
$ws.Cells.Item(14, 1).Value2 = "Municipal Solid Waste (Renewable)"
$ws.Cells.Item(15, 1).Value2 = "Municipal Solid Waste (Non-renewable)"

# Widen columns A and C to fit the new, longer text (bestFit-style resize).
$ws.Columns("A").ColumnWidth = 36.5
$ws.Columns("C").ColumnWidth = 41.6

# Update the active selection to reflect where the editor left off.
$ws.Range("C21").Select() | Out-Null
